# Apply the 2022-10-24 daily crime-data update.
# Updates the 2022 year-to-date totals (column I, plus one 2021 column H
# correction) across the Citywide Totals, By Neighborhood, and each
# individual neighborhood detail sheet, matching the published diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 5956
$ws.Range("I3").Value = 6212
$ws.Range("H4").Value = 1674
$ws.Range("I4").Value = 1426
$ws.Range("I5").Value = 579
$ws.Range("I6").Value = 7038
$ws.Range("H7").Value = 25985
$ws.Range("I7").Value = 21211

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I2").Value = 60
$ws.Range("I7").Value = 242

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I2").Value = 42
$ws.Range("I7").Value = 119

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I3").Value = 227
$ws.Range("I6").Value = 196
$ws.Range("I7").Value = 679

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I3").Value = 140
$ws.Range("I7").Value = 380

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 189
$ws.Range("I3").Value = 301
$ws.Range("I6").Value = 252
$ws.Range("I7").Value = 818

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I2").Value = 66
$ws.Range("I7").Value = 211

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I6").Value = 140
$ws.Range("I7").Value = 494

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I6").Value = 153
$ws.Range("I7").Value = 664
$ws.Range("I8").Value = 1273
$ws.Range("I9").Value = 103
$ws.Range("I10").Value = 150
$ws.Range("I11").Value = 319
$ws.Range("I14").Value = 119
$ws.Range("I15").Value = 245
$ws.Range("I19").Value = 589
$ws.Range("I20").Value = 529
$ws.Range("I22").Value = 56
$ws.Range("I23").Value = 212
$ws.Range("I29").Value = 1313
$ws.Range("I31").Value = 211
$ws.Range("I33").Value = 960
$ws.Range("I34").Value = 99
$ws.Range("I36").Value = 287
$ws.Range("I37").Value = 679
$ws.Range("I42").Value = 729
$ws.Range("I47").Value = 146
$ws.Range("I49").Value = 145
$ws.Range("I51").Value = 245
$ws.Range("I52").Value = 461
$ws.Range("I53").Value = 224
$ws.Range("I54").Value = 435
$ws.Range("I55").Value = 233
$ws.Range("I59").Value = 36
$ws.Range("I60").Value = 115
$ws.Range("H63").Value = 224
$ws.Range("I63").Value = 67
$ws.Range("I65").Value = 494
$ws.Range("I67").Value = 818
$ws.Range("I72").Value = 85
$ws.Range("I73").Value = 195
$ws.Range("I75").Value = 69
$ws.Range("I77").Value = 136
$ws.Range("I78").Value = 287
$ws.Range("I83").Value = 458
$ws.Range("I85").Value = 966
$ws.Range("I88").Value = 193
$ws.Range("I89").Value = 242
$ws.Range("I90").Value = 257
$ws.Range("I93").Value = 122
$ws.Range("I94").Value = 223
$ws.Range("I95").Value = 324
$ws.Range("I97").Value = 180
$ws.Range("I99").Value = 380
$ws.Range("H101").Value = 25985
$ws.Range("I101").Value = 21211

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 157
$ws.Range("I3").Value = 165
$ws.Range("I7").Value = 458

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I3").Value = 118
$ws.Range("I7").Value = 324

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 213
$ws.Range("I3").Value = 365
$ws.Range("I7").Value = 960

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I6").Value = 88
$ws.Range("I7").Value = 145

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I4").Value = 30
$ws.Range("I6").Value = 208
$ws.Range("I7").Value = 435

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I3").Value = 453
$ws.Range("I5").Value = 47
$ws.Range("I6").Value = 364
$ws.Range("I7").Value = 1313

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I3").Value = 181
$ws.Range("I7").Value = 589

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 273
$ws.Range("I3").Value = 372
$ws.Range("I4").Value = 45
$ws.Range("I6").Value = 243
$ws.Range("I7").Value = 966

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I6").Value = 43
$ws.Range("I7").Value = 153

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 182
$ws.Range("I4").Value = 53
$ws.Range("I6").Value = 237
$ws.Range("I7").Value = 729

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I2").Value = 50
$ws.Range("I7").Value = 150

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I3").Value = 71
$ws.Range("I6").Value = 106
$ws.Range("I7").Value = 287

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I3").Value = 74
$ws.Range("I6").Value = 75
$ws.Range("I7").Value = 233

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I4").Value = 12
$ws.Range("I7").Value = 212

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I3").Value = 153
$ws.Range("I6").Value = 184
$ws.Range("I7").Value = 529

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I2").Value = 84
$ws.Range("I6").Value = 91
$ws.Range("I7").Value = 287

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("I2").Value = 33
$ws.Range("I6").Value = 51
$ws.Range("I7").Value = 122

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I3").Value = 163
$ws.Range("I4").Value = 39
$ws.Range("I5").Value = 16
$ws.Range("I7").Value = 461

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("I2").Value = 44
$ws.Range("I7").Value = 99

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I2").Value = 41
$ws.Range("I5").Value = 3
$ws.Range("I7").Value = 223

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I3").Value = 43
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I6").Value = 91
$ws.Range("I7").Value = 245

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I6").Value = 85
$ws.Range("I7").Value = 319

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("I6").Value = 33
$ws.Range("I7").Value = 103

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I3").Value = 61
$ws.Range("I6").Value = 52
$ws.Range("I7").Value = 195

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("I2").Value = 18
$ws.Range("I7").Value = 36

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I6").Value = 117
$ws.Range("I7").Value = 180

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I2").Value = 57
$ws.Range("I7").Value = 193

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 388
$ws.Range("I3").Value = 361
$ws.Range("I4").Value = 76
$ws.Range("I6").Value = 411
$ws.Range("I7").Value = 1273

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("I6").Value = 18
$ws.Range("I7").Value = 69

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I2").Value = 85
$ws.Range("I7").Value = 257

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I6").Value = 103
$ws.Range("I7").Value = 245

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("I2").Value = 40
$ws.Range("I6").Value = 34
$ws.Range("I7").Value = 115

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I4").Value = 20
$ws.Range("I6").Value = 103
$ws.Range("I7").Value = 224

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I6").Value = 16
$ws.Range("I7").Value = 56

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I2").Value = 19
$ws.Range("I7").Value = 85

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("I3").Value = 48
$ws.Range("I7").Value = 136

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 217
$ws.Range("I4").Value = 36
$ws.Range("I6").Value = 174
$ws.Range("I7").Value = 664
